# Fix typo: "sportclub vaak" -> "sportclubs vaak" in the paragraph that
# starts with "Achterliggend worden de gegevens van spelers binnen sportclub ..."
$d = $word.ActiveDocument

$find = $d.Content
$find.Find.Execute("sportclub vaak", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "sportclubs vaak", 2)
